$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.67
$ws.Range("BC2").Value = 140
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("Q3").Value = 2.3
$ws.Range("R3").Value = 1.62
$ws.Range("I4").Value = 4.5
$ws.Range("N4").Value = 7.5
$ws.Range("Z4").Value = 15
$ws.Range("AC4").Value = 7.5
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 19
$ws.Range("G5").Value = 2.55
$ws.Range("H5").Value = 2.82
$ws.Range("I5").Value = 3.1
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 6
$ws.Range("R5").Value = 1.48
$ws.Range("Z5").Value = 26
$ws.Range("AA5").Value = 26
$ws.Range("AI5").Value = 13
$ws.Range("AJ5").Value = 12
$ws.Range("L8").Value = 2.38
$ws.Range("G12").Value = 1.96
$ws.Range("I12").Value = 3.6
$ws.Range("J12").Value = 2.6
$ws.Range("L12").Value = 3.75
$ws.Range("Q12").Value = 1.67
$ws.Range("R12").Value = 2.15
$ws.Range("W12").Value = 9.5
$ws.Range("X12").Value = 11
$ws.Range("Z12").Value = 19
$ws.Range("AI12").Value = 21
$ws.Range("AJ12").Value = 13
$ws.Range("AK12").Value = 41
$ws.Range("AL12").Value = 26
$ws.Range("AQ12").Value = 34
$ws.Range("AX12").Value = 19
$ws.Range("AY12").Value = 23
$ws.Range("I13").Value = 1.62
$ws.Range("G14").Value = 4.1
$ws.Range("I14").Value = 1.62
$ws.Range("J14").Value = 4
$ws.Range("L14").Value = 2.2
$ws.Range("N14").Value = 23
$ws.Range("Z14").Value = 41
$ws.Range("AA14").Value = 26
$ws.Range("AD14").Value = 10
$ws.Range("AS14").Value = 81
$ws.Range("I22").Value = 2.63
$ws.Range("U22").Value = 1.54
$ws.Range("V23").Value = 1.69
$ws.Range("U24").Value = 1.87
$ws.Range("V24").Value = 1.77
$ws.Range("G28").Value = 1.9
$ws.Range("I28").Value = 4.1
$ws.Range("J28").Value = 2.6
$ws.Range("L28").Value = 4.33
$ws.Range("M28").Value = 1.05
$ws.Range("N28").Value = 11
$ws.Range("O28").Value = 1.29
$ws.Range("Q28").Value = 1.92
$ws.Range("R28").Value = 1.82
$ws.Range("U28").Value = 1.8
$ws.Range("V28").Value = 1.91
$ws.Range("X28").Value = 9
$ws.Range("Y28").Value = 8.5
$ws.Range("AA28").Value = 15
$ws.Range("AE28").Value = 15
$ws.Range("AF28").Value = 51
$ws.Range("AH28").Value = 12
$ws.Range("AI28").Value = 21
$ws.Range("AJ28").Value = 15
$ws.Range("AL28").Value = 34
$ws.Range("AM28").Value = 41
$ws.Range("AO28").Value = 10
$ws.Range("AQ28").Value = 34
$ws.Range("AW28").Value = 6
$ws.Range("G30").Value = 2.01
$ws.Range("I30").Value = 3.6
$ws.Range("M30").Value = 1.06
$ws.Range("N30").Value = 10
$ws.Range("O30").Value = 1.3
$ws.Range("Q30").Value = 1.99
$ws.Range("R30").Value = 1.74
$ws.Range("X30").Value = 10
$ws.Range("Z30").Value = 19
$ws.Range("AX30").Value = 19
$ws.Range("M38").Value = 1.05
$ws.Range("O38").Value = 1.41
$ws.Range("P38").Value = 2.62
$ws.Range("M39").Value = 1.03
$ws.Range("O39").Value = 1.25
